$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.302.22"
$ws.Range("E2").Value = "  +1.31%  "

$ws.Range("D3").Value = "1.891.40"
$ws.Range("E3").Value = "  +1.15%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "246.19"
$ws.Range("E5").Value = "  +0.96%  "

$ws.Range("D6").Value = "0.689"
$ws.Range("E6").Value = "  +2.27%  "

$ws.Range("E7").Value = "  +0.02%  "

$ws.Range("D8").Value = "'42.70"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.39%  "

$ws.Range("D9").Value = "0.357"
$ws.Range("E9").Value = "  +4.90%  "

$ws.Range("D10").Value = "56.35"
$ws.Range("E10").Value = "  +10.80%  "

$ws.Range("D11").Value = "0.0747"
$ws.Range("E11").Value = "  +1.95%  "

$ws.Range("E12").Value = "  +1.44%  "

$ws.Range("E13").Value = "  +8.66%  "

$ws.Range("D14").Value = "0.794"
$ws.Range("E14").Value = "  +11.97%  "

$ws.Range("D15").Value = "2.169.59"
$ws.Range("E15").Value = "  +1.45%  "

$ws.Range("D16").Value = "4.99"
$ws.Range("E16").Value = "  +3.68%  "

$ws.Range("D17").Value = "1.931.91"
$ws.Range("E17").Value = "  +3.33%  "

$ws.Range("D18").Value = "35.276.47"
$ws.Range("E18").Value = "  +1.37%  "

$ws.Range("D19").Value = "73.36"
$ws.Range("E19").Value = "  +1.68%  "

$ws.Range("D20").Value = "0.0₃0826"
$ws.Range("E20").Value = "  +2.13%  "

$ws.Range("D21").Value = "243.49"
$ws.Range("E21").Value = "  +0.48%  "

$ws.Range("D22").Value = "'12.90"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +2.83%  "

$ws.Range("D23").Value = "5.22"
$ws.Range("E23").Value = "  +7.11%  "

$ws.Range("E24").Value = "  +8.49%  "

$ws.Range("E25").Value = "  -0.09%  "

$ws.Range("E26").Value = "  +1.29%  "

$ws.Range("D27").Value = "166.76"
$ws.Range("E27").Value = "  +2.22%  "

$ws.Range("D28").Value = "8.52"
$ws.Range("E28").Value = "  +2.49%  "

$ws.Range("D29").Value = "18.27"
$ws.Range("E29").Value = "  +1.51%  "

$ws.Range("E30").Value = "  +1.78%  "

$ws.Range("D31").Value = "0.0604"
$ws.Range("E31").Value = "  +6.33%  "

$ws.Range("D32").Value = "4.33"
$ws.Range("E32").Value = "  +4.22%  "

$ws.Range("D33").Value = "4.22"
$ws.Range("E33").Value = "  +2.56%  "

$ws.Range("D34").Value = "1.87"
$ws.Range("E34").Value = "  +25.26%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("E36").Value = "  -14.71%  "

$ws.Range("D37").Value = "0.849"
$ws.Range("E37").Value = "  +2.67%  "

$ws.Range("E38").Value = "  +1.06%  "

$ws.Range("D39").Value = "0.0727"
$ws.Range("E39").Value = "  +9.42%  "

$ws.Range("E40").Value = "  +7.83%  "

$ws.Range("D41").Value = "98.92"
$ws.Range("E41").Value = "  +1.57%  "

$ws.Range("D42").Value = "16.94"
$ws.Range("E42").Value = "  +0.97%  "

$ws.Range("E43").Value = "  +0.75%  "

$ws.Range("B44").Value = "Gas"
$ws.Range("C44").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D44").Value = "13.65"
$ws.Range("E44").Value = "  +15.74%  "

$ws.Range("B45").Value = "Maker"
$ws.Range("C45").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D45").Value = "1.330.35"
$ws.Range("E45").Value = "  +3.80%  "

$ws.Range("D46").Value = "2.36"
$ws.Range("E46").Value = "  +2.57%  "

$ws.Range("D47").Value = "0.0811"
$ws.Range("E47").Value = "  -1.52%  "

$ws.Range("E48").Value = "  +0.72%  "

$ws.Range("E49").Value = "  +0.30%  "

$ws.Range("E50").Value = "  +1.49%  "

$ws.Range("D51").Value = "42.45"
$ws.Range("E51").Value = "  +0.06%  "
